$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.808.16"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "271.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5287"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3373"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06815"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.95%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7920"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07778"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.841.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.121"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.22%  "

$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.835.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.112.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.673"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.920"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.061"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.399"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.655"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("E28").Value = "  +1.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.324"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.301"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08848"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04949"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.160"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7263"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.874"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.200"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.87%  "

$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5083"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "115.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9333"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.143"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.980"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4419"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1326"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.373"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05936"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.467"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
